$d = $word.ActiveDocument

# --- 1. Amend the opening paragraph: pad trailing spaces, then append the
#        red "(This is a change - Version for main branch)" annotation as
#        three separate runs (matching how Word recorded the original edit). ---
$p1 = $d.Paragraphs(1)
$r = $p1.Range

$r.InsertAfter("  ")

$enDash = [char]0x2013

$seg1Start = $r.End - 1
$r.InsertAfter("(This is a change " + $enDash + " Ve")
$seg1End = $r.End - 1
$d.Range($seg1Start, $seg1End).Font.Color = 255

$seg2Start = $r.End - 1
$r.InsertAfter("rsion for main branch")
$seg2End = $r.End - 1
$d.Range($seg2Start, $seg2End).Font.Color = 255

$seg3Start = $r.End - 1
$r.InsertAfter(")")
$seg3End = $r.End - 1
$d.Range($seg3Start, $seg3End).Font.Color = 255

# --- 2. Remove the trailing "ank God almighty, we are free at last."
#        paragraph entirely (the paragraph that closes out the poem). ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
if ($lastPara.Range.Text -like "*God almighty, we are free at last.*") {
    $lastPara.Range.Delete()
}
